$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 13:22"
$ws.Range("A15").Value = "Araba/Alava"
$ws.Range("A16").Value = "Zaragoza"
$ws.Range("A19").Value = "Salamanca"
$ws.Range("A20").Value = "Valladolid"
$ws.Range("A21").Value = "A Coruña"
$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("A24").Value = "Asturias"
$ws.Range("A25").Value = "Granada"
$ws.Range("A26").Value = "Segovia"
$ws.Range("A27").Value = "Cantabria"
$ws.Range("A33").Value = "Burgos"
$ws.Range("A34").Value = "Cordoba"
$ws.Range("A35").Value = "Jaen"
$ws.Range("A36").Value = "Guadalajara"
$ws.Range("A41").Value = "Avila"
$ws.Range("A42").Value = "Cuenca"
$ws.Range("A45").Value = "Palencia"
$ws.Range("A46").Value = "Lugo"
$ws.Range("A50").Value = "Zamora"
$ws.Range("A51").Value = "Almeria"
$ws.Range("B7").Value = 5392
$ws.Range("C7").Value = 4680
$ws.Range("D7").Value = 4694
$ws.Range("E7").Value = 377
$ws.Range("B15").Value = 2908
$ws.Range("C15").Value = 4680
$ws.Range("D15").Value = 4694
$ws.Range("E15").Value = 245
$ws.Range("B16").Value = 2889
$ws.Range("C16").Value = 597
$ws.Range("D16").Value = 1989
$ws.Range("E16").Value = 303
$ws.Range("B18").Value = 2127
$ws.Range("C18").Value = 451
$ws.Range("D18").Value = 1522
$ws.Range("E18").Value = 154
$ws.Range("B19").Value = 2051
$ws.Range("C19").Value = 542
$ws.Range("D19").Value = 1273
$ws.Range("E19").Value = 236
$ws.Range("B20").Value = 1985
$ws.Range("C20").Value = 724
$ws.Range("D20").Value = 1077
$ws.Range("E20").Value = 184
$ws.Range("B21").Value = 1969
$ws.Range("C21").Value = 333
$ws.Range("D21").Value = 1788
$ws.Range("E21").Value = 67
$ws.Range("B22").Value = 1887
$ws.Range("C22").Value = 200
$ws.Range("D22").Value = 1537
$ws.Range("E22").Value = 150
$ws.Range("B23").Value = 1803
$ws.Range("C23").Value = 4680
$ws.Range("D23").Value = 4694
$ws.Range("E23").Value = 107
$ws.Range("B24").Value = 1799
$ws.Range("C24").Value = 372
$ws.Range("D24").Value = 1306
$ws.Range("E24").Value = 121
$ws.Range("B25").Value = 1686
$ws.Range("C25").Value = 251
$ws.Range("D25").Value = 1287
$ws.Range("E25").Value = 148
$ws.Range("B26").Value = 1672
$ws.Range("C26").Value = 470
$ws.Range("D26").Value = 1070
$ws.Range("E26").Value = 132
$ws.Range("B27").Value = 1659
$ws.Range("C27").Value = 265
$ws.Range("D27").Value = 1292
$ws.Range("E27").Value = 102
$ws.Range("B30").Value = 1513
$ws.Range("C30").Value = 710
$ws.Range("D30").Value = 573
$ws.Range("E30").Value = 230
$ws.Range("B33").Value = 1126
$ws.Range("C33").Value = 470
$ws.Range("D33").Value = 524
$ws.Range("E33").Value = 132
$ws.Range("B34").Value = 1116
$ws.Range("C34").Value = 159
$ws.Range("D34").Value = 909
$ws.Range("E34").Value = 48
$ws.Range("B35").Value = 1055
$ws.Range("C35").Value = 111
$ws.Range("D35").Value = 852
$ws.Range("E35").Value = 92
$ws.Range("B36").Value = 1036
$ws.Range("C36").Value = 1982
$ws.Range("D36").Value = 9650
$ws.Range("E36").Value = 137
$ws.Range("B38").Value = 929
$ws.Range("C38").Value = 144
$ws.Range("D38").Value = 736
$ws.Range("E38").Value = 49
$ws.Range("B40").Value = 898
$ws.Range("C40").Value = 219
$ws.Range("D40").Value = 599
$ws.Range("E40").Value = 80
$ws.Range("B41").Value = 859
$ws.Range("C41").Value = 325
$ws.Range("D41").Value = 441
$ws.Range("E41").Value = 93
$ws.Range("B42").Value = 845
$ws.Range("C42").Value = 1982
$ws.Range("D42").Value = 9650
$ws.Range("E42").Value = 123
$ws.Range("B45").Value = 592
$ws.Range("C45").Value = 162
$ws.Range("D45").Value = 386
$ws.Range("E45").Value = 44
$ws.Range("B46").Value = 586
$ws.Range("C46").Value = 333
$ws.Range("D46").Value = 520
$ws.Range("E46").Value = 11
$ws.Range("B50").Value = 406
$ws.Range("C50").Value = 135
$ws.Range("D50").Value = 222
$ws.Range("E50").Value = 49
$ws.Range("B51").Value = 400
$ws.Range("C51").Value = 73
$ws.Range("D51").Value = 298
$ws.Range("E51").Value = 29
$ws.Range("B52").Value = 310
$ws.Range("C52").Value = 48
$ws.Range("D52").Value = 241
